# Apply "Generate Report for Handback" updates to handback-status.xlsx
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 26b2c152-...md (row 3)
# and 4341e991-...md (row 5) share the same timestamp string; update both.
$wsOverview.Range("G3").Value = "2016-08-13 12:19:57"
$wsOverview.Range("G5").Value = "2016-08-13 12:19:57"

# zh-cn sheet, row for 26b2c152-... (row 3) and 4341e991-... (row 5):
#  - Status changes from "ht" to "mt"
#  - Correspond Handoff Datetime and Correspond Handback DateTime refreshed
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

$wsZhCn.Range("H3").Value = "2016-08-13 12:19:49"
$wsZhCn.Range("H5").Value = "2016-08-13 12:19:49"

$wsZhCn.Range("K3").Value = "2016-08-13 12:20:18"
$wsZhCn.Range("K5").Value = "2016-08-13 12:20:18"

# de-de sheet, row for 26b2c152-... (row 3) and 4341e991-... (row 5):
#  - Correspond Handoff Datetime (shared with Overview's date) refreshed
#  - Correspond Handback DateTime refreshed
$wsDeDe.Range("H3").Value = "2016-08-13 12:19:57"
$wsDeDe.Range("H5").Value = "2016-08-13 12:19:57"

$wsDeDe.Range("K3").Value = "2016-08-13 12:20:27"
$wsDeDe.Range("K5").Value = "2016-08-13 12:20:27"
